# Updated cryptos list on Wed May 22 15:25:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.064.52"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.739.58"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'618.25"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'180.08"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").Value = "3.740.34"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'0.481"
$ws.Range("E12").Value = "  -4.77%  "
$ws.Range("D13").Value = "'40.00"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "4.350.76"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "3.730.95"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "70.063.79"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "'7.56"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'504.50"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'16.38"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").Value = "'9.22"
$ws.Range("E22").Value = "  -3.19%  "
$ws.Range("D23").Value = "'0.720"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "'86.90"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "'12.96"
$ws.Range("E26").Value = "  -4.54%  "
$ws.Range("D27").Value = "'11.13"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'0.0000132"
$ws.Range("E28").Value = "  +7.49%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "'30.45"
$ws.Range("E33").Value = "  -5.51%  "
$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "'0.349"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").Value = "'0.139"
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").Value = "'3.22"
$ws.Range("E40").Value = "  +15.72%  "
$ws.Range("E41").Value = "  -5.13%  "
$ws.Range("D42").Value = "'49.82"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'44.53"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'426.55"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'8.57"
$ws.Range("E45").Value = "  -3.79%  "
$ws.Range("D46").Value = "2.966.26"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("D48").Value = "'27.14"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'136.60"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "'2.48"
$ws.Range("E51").Value = "  -3.06%  "
